# Daily attendance processing - 2026-01-10 21:00:21
#
# The "Recorded By" column (G) lists the users who recorded each session,
# separated by ", ". Re-syncing against the live attendance system
# reorders each cell's contributor list by rotating it one position to
# the left (the first-recorded contributor moves to the end) - unless
# the list already ends with "System", in which case it is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("G" + $row)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ", "
    if ($parts.Count -le 1) { continue }
    if ($parts[$parts.Count - 1] -eq "System") { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
